$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '51.618.68'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '3.028.98'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '384.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.59'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +0.52%  '
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('D13').Value = '3.504.18'
$ws.Range('E13').Value = '  +2.47%  '
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').Value = '3.023.96'
$ws.Range('E16').Value = '  +2.30%  '
$ws.Range('E17').Value = '  -3.13%  '
$ws.Range('E18').Value = '  -11.89%  '
$ws.Range('D19').Value = '51.623.33'
$ws.Range('E19').Value = '  +1.01%  '
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.13%  '
$ws.Range('E26').Value = '  +6.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.54'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.15%  '
$ws.Range('E28').Value = '  +3.60%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.22'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.30%  '
$ws.Range('E31').Value = '  -1.28%  '
$ws.Range('E32').Value = '  -1.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.07'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '34.09'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '50.52'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.33%  '
$ws.Range('E36').Value = '  +3.09%  '
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.299'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.02%  '
$ws.Range('E40').Value = '  +2.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.85'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '127.59'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.40%  '
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.72'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.54'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('E47').Value = '  +2.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.43'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.74%  '
$ws.Range('D49').Value = '2.024.84'
$ws.Range('E49').Value = '  -1.94%  '
$ws.Range('D50').Value = '3.328.85'
$ws.Range('E50').Value = '  +2.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.513'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.12%  '
